# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values (strikeouts) for each game log row.
# These values are raw data (not formulas), regenerated from source stats,
# so they are simply overwritten with the new computed numbers.
$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    8  = 0
    9  = 2
    10 = 2
    11 = 2
    12 = 3
    13 = 0
    14 = 1
    15 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
